# Added New Mac-Address and Document Types
# Append 5 new rows (157-161) to the "master-reg_center_device" data sheet,
# continuing the device_id sequence (3000176 .. 3000180) for regcntr_id 10002,
# keeping the same lang_code/is_active/cr_by/cr_dtimes values used throughout
# the table (eng / TRUE / superadmin / now()).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 157
$startDeviceId = 3000176
$rowCount = 5

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = 10002
    $ws.Cells.Item($r, 2).Value = $startDeviceId + $i
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Match the workbook's view/selection state after the new rows were entered.
$excel.ActiveWindow.ScrollRow = 152
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C158").Select()
